$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: swap the contents of columns B:AD between two rows,
# leaving column A (the sequential row id) untouched.
function Swap-Rows($r1, $r2) {
    $rngA = $ws.Range("B$r1`:AD$r1")
    $rngB = $ws.Range("B$r2`:AD$r2")

    $valsA = $rngA.Value2
    $valsB = $rngB.Value2

    $rngA.Value2 = $valsB
    $rngB.Value2 = $valsA
}

# Simple pairwise swaps
Swap-Rows 224 225
Swap-Rows 231 232
Swap-Rows 237 238
Swap-Rows 249 250

# 3-way rotation for rows 267, 268, 269:
# new267 = old268, new268 = old269, new269 = old267
$r267 = $ws.Range("B267:AD267")
$r268 = $ws.Range("B268:AD268")
$r269 = $ws.Range("B269:AD269")

$v267 = $r267.Value2
$v268 = $r268.Value2
$v269 = $r269.Value2

$r267.Value2 = $v268
$r268.Value2 = $v269
$r269.Value2 = $v267
